# Update the lattice-multiplication exercise table: each table cell holds
# one exercise (problem header, multiplicand digits, separator, and two
# lattice row labels) separated by manual line breaks (vertical-tab chars).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "35 x 80" + [char]11 + "  8    0" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "5|    |"
$t.Cell(1,2).Range.Text = "76 x 81" + [char]11 + "  8    1" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "6|    |"
$t.Cell(1,3).Range.Text = "86 x 33" + [char]11 + "  3    3" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "6|    |"
$t.Cell(2,1).Range.Text = "79 x 76" + [char]11 + "  7    6" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "9|    |"
$t.Cell(2,2).Range.Text = "81 x 96" + [char]11 + "  9    6" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "1|    |"
$t.Cell(2,3).Range.Text = "35 x 65" + [char]11 + "  6    5" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "5|    |"
$t.Cell(3,1).Range.Text = "95 x 50" + [char]11 + "  5    0" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "5|    |"
$t.Cell(3,2).Range.Text = "28 x 69" + [char]11 + "  6    9" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "8|    |"
$t.Cell(3,3).Range.Text = "15 x 89" + [char]11 + "  8    9" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "5|    |"
$t.Cell(4,1).Range.Text = "15 x 89" + [char]11 + "  8    9" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "5|    |"
$t.Cell(4,2).Range.Text = "74 x 72" + [char]11 + "  7    2" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "4|    |"
$t.Cell(4,3).Range.Text = "17 x 72" + [char]11 + "  7    2" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "7|    |"
$t.Cell(5,1).Range.Text = "35 x 26" + [char]11 + "  2    6" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "5|    |"
$t.Cell(5,2).Range.Text = "37 x 80" + [char]11 + "  8    0" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "7|    |"
$t.Cell(5,3).Range.Text = "77 x 76" + [char]11 + "  7    6" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "7|    |"
